# Update the redirects workbook: add a new row for history.html redirects
# (Publication redirects sheet), adjust column B width, and update the
# selected cell to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Publication redirects")

# --- New row 12 data -------------------------------------------------
$ws.Range("A12").Value = "http://ns.electronichealth.net.au/fhir/history.html"
$ws.Range("B12").Value = "https://fhir.digitalhealth.gov.au/dh/history.html"

# Hyperlinks for the new cells (same URL as the displayed text, matching
# the pattern used by the other rows in this table).
$ws.Hyperlinks.Add($ws.Range("A12"), "http://ns.electronichealth.net.au/fhir/history.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B12"), "https://fhir.digitalhealth.gov.au/dh/history.html") | Out-Null

# Copy the formatting (border + hyperlink style) from row 11 down to row 12
# so the new row matches the rest of the table instead of Excel's default
# "just-added-a-hyperlink" style.
$ws.Range("A11:B11").Copy() | Out-Null
$ws.Range("A12:B12").PasteSpecial(-4122) | Out-Null

# --- Column B width ----------------------------------------------------
# Column B is no longer best-fit; it was widened by the author.
$ws.Columns.Item(2).ColumnWidth = 70.8

# --- Final selection -----------------------------------------------------
$ws.Range("B14").Select() | Out-Null
